$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" between "Contact" (row 10) and "Description" (row 11)
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (copy format from row below)
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's content
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
